$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$ws.Range('E2').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F2').Value = 'Ortega Valle Manuel'
$ws.Range('E3').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F3').Value = 'Velasco Sanchez David'
$ws.Range('E4').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F4').Value = 'Hernández Mendoza Delfina'
$ws.Range('E5').Value = 'TEMAS DE FÍSICA'
$ws.Range('F5').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E6').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F6').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E12').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F12').Value = 'Hernández Mendoza Delfina'
$ws.Range('E13').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F13').Value = 'Velasco Sanchez David'
$ws.Range('E14').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F14').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E16').Value = 'TEMAS DE FÍSICA'
$ws.Range('F16').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E17').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F17').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E18').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F18').Value = 'Velasco Sanchez David'
$ws.Range('E19').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F19').Value = 'Hernández Mendoza Delfina'
$ws.Range('E22').Value = 'TEMAS DE FÍSICA'
$ws.Range('F22').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E23').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F23').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E24').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F24').Value = 'Hernández Mendoza Delfina'
$ws.Range('E25').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F25').Value = 'Velasco Sanchez David'
$ws.Range('E30').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F30').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E31').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F31').Value = 'Hernández Mendoza Delfina'
$ws.Range('E35').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F35').Value = 'Velasco Sanchez David'
$ws.Range('E36').Value = 'TEMAS DE FÍSICA'
$ws.Range('F36').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E37').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F37').Value = 'Hernández Mendoza Delfina'
$ws.Range('E38').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F38').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E39').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F39').Value = 'Ortega Valle Manuel'
$ws.Range('E42').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F42').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E43').Value = 'TEMAS DE FÍSICA'
$ws.Range('F43').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E45').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F45').Value = 'Hernández Mendoza Delfina'
$ws.Range('E46').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F46').Value = 'Ortega Valle Manuel'
$ws.Range('E51').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F51').Value = 'Hernández Mendoza Delfina'
$ws.Range('E52').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F52').Value = 'Velasco Sanchez David'
$ws.Range('E54').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F54').Value = 'Velasco Sanchez David'
$ws.Range('E55').Value = 'TEMAS DE FÍSICA'
$ws.Range('F55').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E56').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F56').Value = 'Ortega Valle Manuel'
$ws.Range('E57').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F57').Value = 'Hernández Mendoza Delfina'
